# Re-pull / push all dSF (column F) data for rows 2-34.
# This mirrors an upstream "repull data, push all data, mean calculation"
# commit that recomputed the dSF column values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new dSF (column F) value.
$dsf = @{
    2  = 2
    3  = -5
    4  = -1
    5  = -1
    6  = 0
    7  = 1
    8  = -3
    9  = 0
    10 = -3
    11 = -1
    12 = -4
    13 = -3
    14 = -6
    15 = 0
    16 = 7
    17 = -1
    18 = -1
    19 = 1
    20 = 7
    21 = -4
    22 = 2
    23 = 0
    24 = 6
    25 = 1
    26 = 0
    27 = 4
    28 = -2
    29 = 1
    30 = -3
    31 = -1
    32 = 0
    33 = 2
    34 = -2
}

foreach ($row in $dsf.Keys) {
    $ws.Range("F$row").Value = $dsf[$row]
}
